$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $newBulletText = "1285870 - Marcos Villela Barcza"
    $inserted = $false

    # Preferred approach: build the new paragraph (style "ListBullet" + run
    # text) directly via InsertXML so no incidental session/rsid attributes
    # get stamped onto the new <w:p>.
    try {
        $insertionPoint = $d.Range($target.Range.End, $target.Range.End)

        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' +
               '<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>' + $newBulletText + '</w:t></w:r></w:p>' +
               '<w:p/>' +
               '</w:body></w:document>' +
               '</pkg:xmlData></pkg:part></pkg:package>'

        [void]$insertionPoint.InsertXML($xml)

        # InsertXML leaves a stray trailing empty paragraph behind (it
        # absorbs the paragraph mark that used to separate our new content
        # from whatever followed it). Remove that leftover empty paragraph
        # so the structure matches exactly: Docente(s)... -> new ListBullet
        # paragraph -> Programa resumido.
        foreach ($p in $d.Paragraphs) {
            if ($p.Range.Text -like "$newBulletText*") {
                $stray = $p.Next()
                if ($stray -ne $null -and $stray.Range.Text.Trim().Length -eq 0) {
                    $stray.Range.Delete()
                }
                $inserted = $true
                break
            }
        }
    } catch {
        $inserted = $false
    }

    if (-not $inserted) {
        # Fallback: plain object-model insertion, in case InsertXML is not
        # supported by the runtime in use.
        $target.Range.InsertParagraphAfter()
        $newPara = $target.Next()
        $newPara.Range.Text = $newBulletText
        $newPara.Range.Style = "ListBullet"
    }
}
